# Update the two-digit multiplication problems in the document.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "81×58="; new = "84×12="},
    @{old = "45×71="; new = "14×20="},
    @{old = "39×47="; new = "57×56="},
    @{old = "17×33="; new = "89×28="},
    @{old = "35×13="; new = "43×43="},
    @{old = "65×89="; new = "64×88="},
    @{old = "44×28="; new = "53×50="},
    @{old = "14×41="; new = "70×35="},
    @{old = "47×29="; new = "36×77="},
    @{old = "16×26="; new = "38×37="},
    @{old = "60×28="; new = "52×97="},
    @{old = "52×27="; new = "35×65="},
    @{old = "74×98="; new = "43×70="},
    @{old = "17×50="; new = "58×99="},
    @{old = "52×19="; new = "35×50="},
    @{old = "95×94="; new = "40×41="},
    @{old = "92×17="; new = "60×68="},
    @{old = "51×27="; new = "46×77="},
    @{old = "20×48="; new = "49×23="},
    @{old = "18×86="; new = "79×11="},
    @{old = "80×26="; new = "74×36="},
    @{old = "63×31="; new = "45×38="},
    @{old = "82×86="; new = "48×17="},
    @{old = "12×88="; new = "20×18="},
    @{old = "54×11="; new = "42×13="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
